$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 21 ("DELTA.AIRLINES") into new row 22, preserving formatting
# (fill/border style) as well as the shared string value/type.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value2 = $ws.Range("A21").Value2

# Move the active selection to G23, matching the saved workbook state.
$ws.Range("G23").Select() | Out-Null
